# Replace attendance data: "Anni Hapsah" entries -> "Bashir Rahadi" entries,
# with two additional new rows (7 and 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Bashir Rahadi"
$ws.Range("B2").Value = "2022-07-04 20:42:33"
$ws.Range("C2").ClearContents()

# Row 3
$ws.Range("A3").Value = "Bashir Rahadi"
$ws.Range("B3").Value = "2022-07-05 20:42:33"
$ws.Range("C3").ClearContents()

# Row 4
$ws.Range("A4").Value = "Bashir Rahadi"
$ws.Range("B4").Value = "2022-07-07 09:39:59"
$ws.Range("C4").ClearContents()

# Row 5
$ws.Range("A5").Value = "Bashir Rahadi"
$ws.Range("B5").Value = "2022-07-14 09:12:27"
$ws.Range("C5").Value = "2022-07-14 09:12:53"

# Row 6
$ws.Range("A6").Value = "Bashir Rahadi"
$ws.Range("B6").Value = "2022-07-15 04:51:23"
$ws.Range("C6").Value = "2022-07-15 08:31:03"

# Row 7 (new)
$ws.Range("A7").Value = "Bashir Rahadi"
$ws.Range("B7").Value = "2022-07-17 08:08:40"

# Row 8 (new)
$ws.Range("A8").Value = "Bashir Rahadi"
$ws.Range("B8").Value = "2022-07-17 08:08:40"
